# Added for RB #3:
#  - Update the shared "Период" value used by sheet "Скидка за объем закупа"
#  - Add new worksheet "Скидка  на группы товаров" with its discount-by-group table

$wb = $excel.ActiveWorkbook

# 1) Update the period text (shared by the existing "Скидка за объем закупа" sheet too)
$wsVolume = $wb.Worksheets.Item("Скидка за объем закупа")
$wsVolume.Range("A2").Value = "01.03.2022-01.01.2022"

# 2) Add the new worksheet at the end of the workbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Скидка  на группы товаров"

# Grab the header/total style already used on the existing sheets (fill + thin border)
$styleSrc = $wsVolume.Range("A1").Style

# Header row
$ws.Range("A1").Value = "Период"
$ws.Range("B1").Value = "Номер договора/ДС"
$ws.Range("C1").Value = "Тип скидки"
$ws.Range("D1").Value = "Код товара"
$ws.Range("E1").Value = "План закупа"
$ws.Range("F1").Value = "Скидка %"
$ws.Range("G1").Value = "Сумма скидки"
$ws.Range("A1:G1").Style = $styleSrc

# Data rows
$ws.Range("A2").Value = "01.03.2022-01.01.2022"
$ws.Range("B2").Value = "2500800DLR"
$ws.Range("C2").Value = "Скидка  на группы товаров"
$ws.Range("D2").Value = "00000064865"
$ws.Range("E2").Value = 500
$ws.Range("F2").Value = 12
$ws.Range("G2").Value = 295.44

$ws.Range("A3").Value = "01.03.2022-01.01.2022"
$ws.Range("B3").Value = "2500800DLR"
$ws.Range("C3").Value = "Скидка  на группы товаров"
$ws.Range("D3").Value = "00000045698"
$ws.Range("E3").Value = 1000
$ws.Range("F3").Value = 14
$ws.Range("G3").Value = 0

$ws.Range("A4").Value = "01.03.2022-01.01.2022"
$ws.Range("B4").Value = "2500800DLR"
$ws.Range("C4").Value = "Скидка  на группы товаров"
$ws.Range("D4").Value = "00000053058"
$ws.Range("E4").Value = 1500
$ws.Range("F4").Value = 14
$ws.Range("G4").Value = 0

# Totals row
$ws.Range("F5").Value = "Итог:"
$ws.Range("G5").Value = 295.44
$ws.Range("F5:G5").Style = $styleSrc
